$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.488.24'
$ws.Range('E2').Value = '  -3.25%  '
$ws.Range('D3').Value = '1.862.13'
$ws.Range('E3').Value = '  -4.31%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.007'
$ws.Range('E4').Value = '  -0.98%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '320.99'
$ws.Range('E5').Value = '  +0.02%  '
$ws.Range('E6').Value = '  -0.78%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4500'
$ws.Range('E7').Value = '  -5.35%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3868'
$ws.Range('E8').Value = '  -3.72%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '47.74'
$ws.Range('E9').Value = '  -10.90%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07994'
$ws.Range('E10').Value = '  -5.82%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.019'
$ws.Range('E11').Value = '  -3.66%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '21.51'
$ws.Range('E12').Value = '  -2.21%  '
$ws.Range('D13').Value = '1.917.31'
$ws.Range('E13').Value = '  -3.41%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.896'
$ws.Range('E14').Value = '  -4.66%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.153'
$ws.Range('E15').Value = '  -5.54%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.009'
$ws.Range('E16').Value = '  -0.79%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001036'
$ws.Range('E17').Value = '  -3.60%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '86.09'
$ws.Range('E18').Value = '  -2.88%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06547'
$ws.Range('E19').Value = '  -1.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.15'
$ws.Range('E20').Value = '  -8.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.005'
$ws.Range('E21').Value = '  -0.98%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.518'
$ws.Range('E22').Value = '  -5.05%  '
$ws.Range('D23').Value = '27.543.73'
$ws.Range('E23').Value = '  -3.23%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.87'
$ws.Range('E24').Value = '  -5.33%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.280'
$ws.Range('E25').Value = '  -0.57%  '
$ws.Range('D26').Value = '2.126.27'
$ws.Range('E26').Value = '  -4.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '151.29'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.47'
$ws.Range('E28').Value = '  -3.22%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.514'
$ws.Range('E29').Value = '  -6.78%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.031'
$ws.Range('E30').Value = '  -5.59%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '121.44'
$ws.Range('E31').Value = '  -1.56%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09401'
$ws.Range('E32').Value = '  -1.46%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.464'
$ws.Range('E33').Value = '  +2.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9272'
$ws.Range('E34').Value = '  -6.44%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.634'
$ws.Range('E35').Value = '  -1.10%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.282'
$ws.Range('E36').Value = '  -5.66%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.226'
$ws.Range('E37').Value = '  -2.46%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02228'
$ws.Range('E38').Value = '  -4.36%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05994'
$ws.Range('E39').Value = '  -3.36%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.365'
$ws.Range('E40').Value = '  -3.92%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.006'
$ws.Range('E41').Value = '  -0.79%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5946'
$ws.Range('E42').Value = '  -4.44%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1866'
$ws.Range('E43').Value = '  -2.80%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.31'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.280'
$ws.Range('E45').Value = '  -3.54%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5675'
$ws.Range('E46').Value = '  -4.47%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '12.20'
$ws.Range('E47').Value = '  -6.11%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.939'
$ws.Range('E48').Value = '  -5.71%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.371'
$ws.Range('E49').Value = '  -0.91%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06857'
$ws.Range('E50').Value = '  +0.72%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.006'
$ws.Range('E51').Value = '  -0.96%  '
